$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "Keterangan"/Kolom-header template data (columns A:F and M:N)
$ws.Cells.Clear()

# New header row
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Nama Pegawai"
$ws.Range("C1").Value = "Divisi"
$ws.Range("D1").Value = "NIP"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Pangkat"

# New data row - imported employee record
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "mutia"
$ws.Range("C2").Value = "Produksi"
$ws.Range("D2").Value = 12345
$ws.Range("E2").Value = "mutia@bps.go.id"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:mutia@bps.go.id")
$ws.Range("F2").Value = "sekretaris"

# Selection ends on J9, matching the saved workbook view state
$ws.Range("J9").Select()
